$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "MORENO CANCHANYA, ROSMERY"
$ws.Cells.Item(2, 4).Value = "MORENO CANCHANYA, ROSMERY"


$ws.Cells.Item(4, 4).Value = "HUAMAN HUAMANI, ALEXIS JAVIER"

$ws.Cells.Item(5, 1).Value = "HUAMAN HUAMANI, ALEXIS JAVIER"
$ws.Cells.Item(5, 4).Value = "YOVERA ROBLES, VICTOR EDUARDO"

$ws.Cells.Item(6, 4).Value = "MUÑOZ SOTOMAYOR, MIRIAN RAQUEL"

$ws.Cells.Item(7, 1).Value = "MUÑOZ SOTOMAYOR, MIRIAN RAQUEL"
$ws.Cells.Item(7, 5).Value = "12:45"
$ws.Cells.Item(7, 6).Value = "09:00"

$ws.Cells.Item(8, 2).Value = "09:00"
$ws.Cells.Item(8, 3).Value = "12:45"
$ws.Cells.Item(8, 4).Value = "QUIQUIA MALLQUI, CYNTHIA ANGELLINE"
$ws.Cells.Item(8, 5).Value = "13:45"
$ws.Cells.Item(8, 6).Value = "10:00"

$ws.Cells.Item(9, 1).Value = "QUIQUIA MALLQUI, CYNTHIA ANGELLINE"
$ws.Cells.Item(9, 2).Value = "10:00"
$ws.Cells.Item(9, 3).Value = "13:45"
$ws.Cells.Item(9, 4).Value = "POBLETE SAIRE, FIORELLA ESTHER"

$ws.Cells.Item(10, 1).Value = "POBLETE SAIRE, FIORELLA ESTHER"
$ws.Cells.Item(10, 4).Value = "GOMEZ ALBINO, IDALIA GIMENA"

$ws.Cells.Item(11, 1).Value = "AGUILAR SCHLAEFLI, STEPHANIE XIMENA"
$ws.Cells.Item(11, 4).Value = "AGUILAR SCHLAEFLI, STEPHANIE XIMENA"

$ws.Cells.Item(12, 1).Value = "GOMEZ ALBINO, IDALIA GIMENA"
$ws.Cells.Item(12, 4).Value = "RUIZ SANTOS, CIELO CRISTHINA"
$ws.Cells.Item(12, 5).Value = "14:00"
$ws.Cells.Item(12, 6).Value = "10:15"

$ws.Cells.Item(13, 1).Value = "RUIZ SANTOS, CIELO CRISTHINA"
$ws.Cells.Item(13, 2).Value = "10:15"
$ws.Cells.Item(13, 3).Value = "14:00"
$ws.Cells.Item(13, 4).Value = "CUSI QUISPE, ANDREA ESTEFANY"

$ws.Cells.Item(14, 1).Value = "CUSI QUISPE, ANDREA ESTEFANY"
$ws.Cells.Item(14, 4).Value = "HUAYANAY VELASCO, ATHINA"

$ws.Cells.Item(15, 1).Value = "HUAYANAY VELASCO, ATHINA"
$ws.Cells.Item(15, 4).Value = "MONTEZUMA DEJO, EVELYN BRUNELLA"

$ws.Cells.Item(16, 1).Value = "MONTEZUMA DEJO, EVELYN BRUNELLA"
$ws.Cells.Item(16, 4).Value = "QUISPE MONDRAGÓN, JUAN ALFONSO"

$ws.Cells.Item(17, 1).Value = "QUISPE MONDRAGÓN, JUAN ALFONSO"
$ws.Cells.Item(17, 4).Value = "AYALA MORA, CECILIA ROSARIO"

$ws.Cells.Item(18, 1).Value = "AYALA MORA, CECILIA ROSARIO"
$ws.Cells.Item(18, 4).Value = "VARGAS CASTRO, LOANA VICTORIA"

$ws.Cells.Item(19, 2).Value = "12:00"
$ws.Cells.Item(19, 3).Value = "21:00"

$ws.Cells.Item(20, 1).Value = "ZEVALLOS ZANCA, VERONICA LUZ"
$ws.Cells.Item(20, 3).Value = "23:00"
$ws.Cells.Item(20, 4).Value = "FLORES PAREDES, LOURDES"
$ws.Cells.Item(20, 6).Value = "15:30"

$ws.Cells.Item(21, 1).Value = "TITO LAURA, NANCY FIORELLA"
$ws.Cells.Item(21, 2).Value = "14:00"
$ws.Cells.Item(21, 3).Value = "23:00"
$ws.Cells.Item(21, 4).Value = "TORRES RAZURI, JESUS GUSTAVO SANTIAGO"

$ws.Cells.Item(22, 1).Value = "VARGAS CASTRO, LOANA VICTORIA"
$ws.Cells.Item(22, 2).Value = "14:00"
$ws.Cells.Item(22, 3).Value = "17:45"
$ws.Cells.Item(22, 4).Value = "ALVITE CORNEJO, ANGIE LUCERO"

$ws.Cells.Item(23, 1).Value = "TORRES RAZURI, JESUS GUSTAVO SANTIAGO"
$ws.Cells.Item(23, 2).Value = "14:30"
$ws.Cells.Item(23, 3).Value = "18:15"
$ws.Cells.Item(23, 4).Value = "SUAREZ JARA, YENNIFER YUSSARA"

$ws.Cells.Item(24, 1).Value = "ALVITE CORNEJO, ANGIE LUCERO"
$ws.Cells.Item(24, 2).Value = "14:45"
$ws.Cells.Item(24, 3).Value = "18:30"
$ws.Cells.Item(24, 4).Value = "ARIAS MACHACUAY, SADELITH SORAGGI"
$ws.Cells.Item(24, 5).Value = "18:45"
$ws.Cells.Item(24, 6).Value = "15:00"

$ws.Cells.Item(25, 1).Value = "SUAREZ JARA, YENNIFER YUSSARA"
$ws.Cells.Item(25, 2).Value = "15:00"
$ws.Cells.Item(25, 3).Value = "18:45"
$ws.Cells.Item(25, 4).Value = "AYALA TAPIA, DARCIE SOL"
$ws.Cells.Item(25, 5).Value = "20:15"
$ws.Cells.Item(25, 6).Value = "16:30"

$ws.Cells.Item(26, 1).Value = "ARIAS MACHACUAY, SADELITH SORAGGI"
$ws.Cells.Item(26, 2).Value = "15:00"
$ws.Cells.Item(26, 3).Value = "18:45"
$ws.Cells.Item(26, 4).Value = "CARDENAS RICAPA, FABRIZIO ESTEBAN"

$ws.Cells.Item(27, 1).Value = "FLORES PAREDES, LOURDES"
$ws.Cells.Item(27, 2).Value = "15:30"
$ws.Cells.Item(27, 3).Value = "18:15"
$ws.Cells.Item(27, 4).Value = "ILDEFONSO MOTTA, JHOSSEP ANGELO"
$ws.Cells.Item(27, 5).Value = "20:30"
$ws.Cells.Item(27, 6).Value = "16:45"

$ws.Cells.Item(28, 1).Value = "AYALA TAPIA, DARCIE SOL"
$ws.Cells.Item(28, 4).Value = "YANAC DAVILA, GERALD RONNY"
$ws.Cells.Item(28, 5).Value = "20:45"
$ws.Cells.Item(28, 6).Value = "17:00"

$ws.Cells.Item(29, 1).Value = "CARDENAS RICAPA, FABRIZIO ESTEBAN"
$ws.Cells.Item(29, 2).Value = "16:30"
$ws.Cells.Item(29, 3).Value = "20:15"
$ws.Cells.Item(29, 4).Value = "SOTELO GONZALES, CAMILA SOFÍA"
$ws.Cells.Item(29, 5).Value = "21:00"
$ws.Cells.Item(29, 6).Value = "12:00"

$ws.Cells.Item(30, 1).Value = "ILDEFONSO MOTTA, JHOSSEP ANGELO"
$ws.Cells.Item(30, 2).Value = "16:45"
$ws.Cells.Item(30, 3).Value = "20:30"
$ws.Cells.Item(30, 4).Value = "SALAS VILLANUEVA, JAMILA DASHA"

$ws.Cells.Item(31, 1).Value = "YANAC DAVILA, GERALD RONNY"
$ws.Cells.Item(31, 2).Value = "17:00"
$ws.Cells.Item(31, 3).Value = "20:45"
$ws.Cells.Item(31, 4).Value = "CAPCHA YARANGO, DAVID"

$ws.Cells.Item(32, 1).Value = "SALAS VILLANUEVA, JAMILA DASHA"
$ws.Cells.Item(32, 2).Value = "17:30"
$ws.Cells.Item(32, 3).Value = "21:15"
$ws.Cells.Item(32, 4).Value = "VILCHEZ CUBA, JACK ANTHONY"
$ws.Cells.Item(32, 6).Value = "17:45"

$ws.Cells.Item(33, 1).Value = "VILCHEZ CUBA, JACK ANTHONY"
$ws.Cells.Item(33, 2).Value = "17:45"
$ws.Cells.Item(33, 3).Value = "21:30"
$ws.Cells.Item(33, 4).Value = "INGA DELGADO, CARLOS DANIEL"

$ws.Cells.Item(34, 1).Value = "CAPCHA YARANGO, DAVID"
$ws.Cells.Item(34, 2).Value = "17:45"
$ws.Cells.Item(34, 3).Value = "21:30"
$ws.Cells.Item(34, 4).Value = "MARTICORENA LOPEZ, DAVID CARLOS"

$ws.Cells.Item(35, 1).Value = "MARTICORENA LOPEZ, DAVID CARLOS"
$ws.Cells.Item(35, 4).Value = "BONILLA SÁNCHEZ, RAÚL FERNANDO"
$ws.Cells.Item(35, 5).Value = "22:00"
$ws.Cells.Item(35, 6).Value = "18:15"

$ws.Cells.Item(36, 1).Value = "INGA DELGADO, CARLOS DANIEL"
$ws.Cells.Item(36, 2).Value = "18:00"
$ws.Cells.Item(36, 3).Value = "21:45"
$ws.Cells.Item(36, 4).Value = "MENDOZA CRUZ, LILIANA LILIANA"

$ws.Cells.Item(37, 1).Value = "MENDOZA CRUZ, LILIANA LILIANA"
$ws.Cells.Item(37, 4).Value = "BRENIS LÁRTIGA, SEBASTIÁN"
$ws.Cells.Item(37, 5).Value = "22:15"
$ws.Cells.Item(37, 6).Value = "18:30"

$ws.Cells.Item(38, 1).Value = "BONILLA SÁNCHEZ, RAÚL FERNANDO"
$ws.Cells.Item(38, 4).Value = "MEDINA MARCELO, NAOMI ARIADNA"
$ws.Cells.Item(38, 5).Value = "22:30"
$ws.Cells.Item(38, 6).Value = "18:45"

$ws.Cells.Item(39, 4).Value = "CORDOVA MONTALVO, MELANY KARINA"
$ws.Cells.Item(39, 5).Value = "22:45"
$ws.Cells.Item(39, 6).Value = "19:00"

$ws.Cells.Item(40, 1).Value = "MEDINA MARCELO, NAOMI ARIADNA"
$ws.Cells.Item(40, 2).Value = "18:45"
$ws.Cells.Item(40, 3).Value = "22:30"
$ws.Cells.Item(40, 4).Value = "LA ROSA EUSEBIO, SHADIA SHAMIRA"

$ws.Cells.Item(41, 1).Value = "CORDOVA MONTALVO, MELANY KARINA"
$ws.Cells.Item(41, 4).Value = "TITO LAURA, NANCY FIORELLA"
$ws.Cells.Item(41, 5).Value = "23:00"
$ws.Cells.Item(41, 6).Value = "14:00"

$ws.Cells.Item(42, 1).Value = "LA ROSA EUSEBIO, SHADIA SHAMIRA"
$ws.Cells.Item(42, 2).Value = "19:00"
$ws.Cells.Item(42, 3).Value = "22:45"
$ws.Cells.Item(42, 4).Value = "ZEVALLOS ZANCA, VERONICA LUZ"
$ws.Cells.Item(42, 5).Value = "23:00"
$ws.Cells.Item(42, 6).Value = "14:00"
